$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset dropped two teams - "BORGES CLIMA FUT F.C" (previously row 7)
# and "BORGES ITAQUI F.C." (previously row 8). Remove the existing
# hyperlinks first (row deletion does not renumber/refresh the hyperlink
# locations on its own), delete the two rows - which shifts every row below
# up by two - and then rebuild the "Link do Time" hyperlinks so they match
# the new layout.
$ws.Hyperlinks.Delete()
$ws.Rows("7:8").Delete()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $id = $ws.Cells.Item($r, 2).Value()
    $location = "!/time/" + $id
    $cell = $ws.Cells.Item($r, 3)
    $ws.Hyperlinks.Add($cell, "https://cartola.globo.com/", $location) | Out-Null
}

# Hyperlinks.Add() re-applies the built-in "Hyperlink" cell style, which can
# leave a near-duplicate style entry behind; explicitly re-stamping the
# "Hyperlink" style on the whole column keeps every cell pointed at the
# original style index.
$ws.Range("C2:C" + $lastRow).Style = "Hyperlink"
